$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (pushes the existing rows 85-184 down to 86-185,
# extending the used range from A1:R184 to A1:R185 — a new week of price data for
# "Vega Monumental Concepción - Pepino ensalada" was added at the top of the table).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record.
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value = "Bíobío"
$ws.Cells.Item(85, 4).Value = 44966
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112043
$ws.Cells.Item(85, 7).Value = "Pepino ensalada"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 100
$ws.Cells.Item(85, 11).Value = 8000
$ws.Cells.Item(85, 12).Value = 9000
$ws.Cells.Item(85, 13).Value = 8500
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 142
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
